$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# The table's current last row holds the "naturalenv_300m" data. We need to:
#  1. Insert a brand-new row just before it, carrying that same
#     "naturalenv_300m" data forward (this becomes the new second-to-last row).
#  2. Turn the former last row into the new "Blood glucose" row.
$lastRow = $t.Rows.Item($t.Rows.Count)
$t.Rows.Add($lastRow) | Out-Null

# Re-fetch by index: inserting shifts the pre-existing last row down one slot,
# so it is now the true last row again, and the freshly-inserted blank row
# sits just above it.
$newRow = $t.Rows.Item($t.Rows.Count - 1)
$lastRow = $t.Rows.Item($t.Rows.Count)

$newRow.Cells.Item(1).Range.Text = "naturalenv_300m, Median (Q1, Q3)"
$newRow.Cells.Item(2).Range.Text = "21.2 (7.4, 42.4)"
$newRow.Cells.Item(3).Range.Text = "19.4 (5.8, 41.1)"

$lastRow.Cells.Item(1).Range.Text = "Blood glucose, Median (Q1, Q3)"
$lastRow.Cells.Item(2).Range.Text = "36.0 (33.4, 38.5)"
$lastRow.Cells.Item(3).Range.Text = "34.7 (32.4, 37.2)"

Write-Output "rows now: $($t.Rows.Count)"
